$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns stay text so numeric-looking
# strings like "1.000" or "0.9999" are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.635.84"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "1.869.20"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "235.24"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "0.4709"
$ws.Range("E7").Value = "  -1.30%  "
$ws.Range("D8").Value = "0.2770"
$ws.Range("E8").Value = "  +0.69%  "
$ws.Range("D9").Value = "0.06387"
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").Value = "17.98"
$ws.Range("E10").Value = "  +11.01%  "
$ws.Range("D11").Value = "1.866.68"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "0.07459"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "4.991"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "85.32"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "0.6371"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").Value = "30.594.68"
$ws.Range("E16").Value = "  +0.97%  "
$ws.Range("D17").Value = "241.27"
$ws.Range("E17").Value = "  +2.70%  "
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").Value = "12.90"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").Value = "0.000007395"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "5.001"
$ws.Range("E22").Value = "  -1.78%  "
$ws.Range("D23").Value = "6.065"
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("D24").Value = "9.417"
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("D25").Value = "165.80"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("D26").Value = "18.24"
$ws.Range("E26").Value = "  +2.00%  "
$ws.Range("D27").Value = "1.896"
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("D28").Value = "0.1023"
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("D29").Value = "1.383"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "4.104"
$ws.Range("E30").Value = "  -2.52%  "
$ws.Range("D31").Value = "3.876"
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("D32").Value = "0.04936"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").Value = "1.154"
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("D34").Value = "0.7124"
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("D35").Value = "2.709"
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("D36").Value = "0.01911"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").Value = "2.696"
$ws.Range("E37").Value = "  +2.31%  "
$ws.Range("D38").Value = "0.8829"
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("D39").Value = "1.998"
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("D40").Value = "105.92"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").Value = "0.9999"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "0.4118"
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("D43").Value = "5.567"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "7.407"
$ws.Range("E44").Value = "  +4.77%  "
$ws.Range("D45").Value = "62.88"
$ws.Range("E45").Value = "  +2.78%  "
$ws.Range("D46").Value = "0.1233"
$ws.Range("E46").Value = "  +2.22%  "
$ws.Range("D47").Value = "33.75"
$ws.Range("E47").Value = "  +2.16%  "
$ws.Range("D48").Value = "8.637"
$ws.Range("E48").Value = "  -1.16%  "
$ws.Range("D49").Value = "0.05574"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").Value = "1.382"
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("D51").Value = "0.3717"
$ws.Range("E51").Value = "  +0.58%  "
